# Scheduled market-data refresh: update cached price/profit figures for the
# Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Values are
# static snapshots (no formulas in this workbook), so each changed cell is
# written directly with the freshly-pulled figure.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 315.7
$ws.Range("J2").Value = 400
$ws.Range("L2").Value = 400
$ws.Range("N2").Value = -626

$ws.Range("H17").Value = 1816.6818
$ws.Range("J17").Value = 1816.6818
$ws.Range("L17").Value = 5450.0454
$ws.Range("N17").Value = -5786.0454

$ws.Range("H43").Value = 2139
$ws.Range("I43").Value = 683.3333
$ws.Range("K43").Value = 683.3333
$ws.Range("M43").Value = -614.3333

$ws.Range("H74").Value = 6320.758
$ws.Range("I74").Value = 6502.067
$ws.Range("J74").Value = 6169.6665
$ws.Range("K74").Value = 6502.067
$ws.Range("L74").Value = 6169.6665
$ws.Range("M74").Value = -5566.067
$ws.Range("N74").Value = -8041.6665

$ws.Range("H77").Value = 6320.758
$ws.Range("I77").Value = 6502.067
$ws.Range("J77").Value = 6169.6665
$ws.Range("K77").Value = 32510.335
$ws.Range("L77").Value = 30848.3325
$ws.Range("M77").Value = -27830.335
$ws.Range("N77").Value = -40208.3325

$ws.Range("H86").Value = 4008.8
$ws.Range("I86").Value = 2303
$ws.Range("K86").Value = 2303
$ws.Range("M86").Value = -1180

$ws.Range("H89").Value = 4008.8
$ws.Range("I89").Value = 2303
$ws.Range("K89").Value = 11515
$ws.Range("M89").Value = -5899

$ws.Range("H132").Value = 4379.9443
$ws.Range("I132").Value = 4573.303
$ws.Range("J132").Value = 2253
$ws.Range("K132").Value = 13719.909
$ws.Range("L132").Value = 6759
$ws.Range("M132").Value = -11189.909
$ws.Range("N132").Value = -11819

$ws.Range("H136").Value = 54999.8
$ws.Range("J136").Value = 54999.8
$ws.Range("L136").Value = 54999.8
$ws.Range("N136").Value = -65199.8

$ws.Range("H137").Value = 2465.1667
$ws.Range("I137").Value = 2572.75
$ws.Range("K137").Value = 7718.25
$ws.Range("M137").Value = -5168.25

$ws.Range("H138").Value = 8038.8716
$ws.Range("J138").Value = 6671.7095
$ws.Range("L138").Value = 20015.1285
$ws.Range("N138").Value = -30295.1285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3740.4878
$ws.Range("I61").Value = 2585.9429
$ws.Range("K61").Value = 2585.9429
$ws.Range("M61").Value = -2373.9429

$ws.Range("H74").Value = 6365.174
$ws.Range("I74").Value = 8612.143
$ws.Range("J74").Value = 2869.889
$ws.Range("K74").Value = 8612.143
$ws.Range("L74").Value = 2869.889
$ws.Range("M74").Value = -7738.143
$ws.Range("N74").Value = -4617.889

$ws.Range("H77").Value = 6365.174
$ws.Range("I77").Value = 8612.143
$ws.Range("J77").Value = 2869.889
$ws.Range("K77").Value = 43060.715
$ws.Range("L77").Value = 14349.445
$ws.Range("M77").Value = -38692.715
$ws.Range("N77").Value = -23085.445

$ws.Range("H102").Value = 6378.1055
$ws.Range("I102").Value = 10700
$ws.Range("J102").Value = 1576
$ws.Range("K102").Value = 10700
$ws.Range("L102").Value = 1576
$ws.Range("M102").Value = -9078
$ws.Range("N102").Value = -4820

$ws.Range("H122").Value = 9863.65
$ws.Range("I122").Value = 9863.65
$ws.Range("K122").Value = 29590.95
$ws.Range("M122").Value = -27140.95

$ws.Range("H132").Value = 14288943
$ws.Range("I132").Value = 18184376
$ws.Range("K132").Value = 54553128
$ws.Range("M132").Value = -54550598

$ws.Range("H136").Value = 3740.4878
$ws.Range("I136").Value = 2585.9429
$ws.Range("K136").Value = 7757.8287
$ws.Range("M136").Value = -5207.8287

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7036.5903
$ws.Range("I134").Value = 6114.104
$ws.Range("K134").Value = 18342.312
$ws.Range("M134").Value = -15807.312

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 9941
$ws.Range("I122").Value = 9941
$ws.Range("K122").Value = 29823
$ws.Range("M122").Value = -27373

$ws.Range("H132").Value = 58829324
$ws.Range("J132").Value = 6357.25
$ws.Range("L132").Value = 19071.75
$ws.Range("N132").Value = -24131.75

$ws.Range("H141").Value = 499607.34
$ws.Range("I141").Value = 69648.5
$ws.Range("J141").Value = 585599.1
$ws.Range("K141").Value = 69648.5
$ws.Range("L141").Value = 585599.1
$ws.Range("M141").Value = -64468.5
$ws.Range("N141").Value = -595959.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 2657.1428
$ws.Range("I92").Value = 1199
$ws.Range("J92").Value = 2900.1667
$ws.Range("K92").Value = 3597
$ws.Range("L92").Value = 8700.500100000001
$ws.Range("M92").Value = -2349
$ws.Range("N92").Value = -11196.5001

$ws.Range("H98").Value = 35715588
$ws.Range("J98").Value = 83335624
$ws.Range("L98").Value = 250006872
$ws.Range("N98").Value = -250009868

$ws.Range("H131").Value = 40517744
$ws.Range("J131").Value = 37175824
$ws.Range("L131").Value = 111527472
$ws.Range("N131").Value = -111537552

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3000
$ws.Range("I80").Value = 3000
$ws.Range("K80").Value = 3000
$ws.Range("M80").Value = -2002

$ws.Range("H83").Value = 3000
$ws.Range("I83").Value = 3000
$ws.Range("K83").Value = 15000
$ws.Range("M83").Value = -10008

$ws.Range("H102").Value = 17107.223
$ws.Range("I102").Value = 18853.5
$ws.Range("K102").Value = 18853.5
$ws.Range("M102").Value = -17231.5

$ws.Range("H122").Value = 93986.55
$ws.Range("I122").Value = 102885.2
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 308655.6
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -306205.6
$ws.Range("N122").Value = -19900

$ws.Range("H132").Value = 18199546
$ws.Range("I132").Value = 40007544
$ws.Range("J132").Value = 26214.5
$ws.Range("K132").Value = 120022632
$ws.Range("L132").Value = 78643.5
$ws.Range("M132").Value = -120020102
$ws.Range("N132").Value = -83703.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2634.923
$ws.Range("I7").Value = 2644.3333
$ws.Range("K7").Value = 2644.3333
$ws.Range("M7").Value = -2532.3333

$ws.Range("H22").Value = 1056.3914
$ws.Range("I22").Value = 865.3125
$ws.Range("K22").Value = 865.3125
$ws.Range("M22").Value = -570.3125

$ws.Range("H27").Value = 1056.3914
$ws.Range("I27").Value = 865.3125
$ws.Range("K27").Value = 865.3125
$ws.Range("M27").Value = -758.3125

$ws.Range("H56").Value = 34999
$ws.Range("I56").Value = 19999
$ws.Range("K56").Value = 19999
$ws.Range("M56").Value = -19308

$ws.Range("H68").Value = 2040.8462
$ws.Range("I68").Value = 1957.6364
$ws.Range("K68").Value = 1957.6364
$ws.Range("M68").Value = -1208.6364

$ws.Range("H71").Value = 2040.8462
$ws.Range("I71").Value = 1957.6364
$ws.Range("K71").Value = 9788.182000000001
$ws.Range("M71").Value = -6044.182000000001

$ws.Range("H82").Value = 5608.222
$ws.Range("I82").Value = 995.8333
$ws.Range("J82").Value = 14833
$ws.Range("K82").Value = 995.8333
$ws.Range("L82").Value = 14833
$ws.Range("M82").Value = -634.8333
$ws.Range("N82").Value = -15555

$ws.Range("H85").Value = 5608.222
$ws.Range("I85").Value = 995.8333
$ws.Range("J85").Value = 14833
$ws.Range("K85").Value = 995.8333
$ws.Range("L85").Value = 14833
$ws.Range("M85").Value = 252.1667
$ws.Range("N85").Value = -17329

$ws.Range("H103").Value = 51664.332
$ws.Range("J103").Value = 51664.332
$ws.Range("L103").Value = 51664.332
$ws.Range("N103").Value = -54008.332

$ws.Range("H126").Value = 2634.923
$ws.Range("I126").Value = 2644.3333
$ws.Range("K126").Value = 7932.999899999999
$ws.Range("M126").Value = -5462.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 16625.5
$ws.Range("J4").Value = 2500
$ws.Range("L4").Value = 2500
$ws.Range("N4").Value = -2726

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws.Range("H122").Value = 102516.4
$ws.Range("I122").Value = 113684.89
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 341054.67
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -338604.67
$ws.Range("N122").Value = -10900

$ws.Range("H124").Value = 33749.92
$ws.Range("J124").Value = 32395.75
$ws.Range("L124").Value = 32395.75
$ws.Range("N124").Value = -42215.75

$ws.Range("H132").Value = 19080.5
$ws.Range("I132").Value = 22197.6
$ws.Range("K132").Value = 66592.79999999999
$ws.Range("M132").Value = -64062.79999999999

$ws.Range("H136").Value = 8714.293
$ws.Range("I136").Value = 8290.102000000001
$ws.Range("K136").Value = 24870.306
$ws.Range("M136").Value = -22320.306
